# v1.3.6 hotfix for 1.3.5
# Append additional filtered chat rows (shared strings + rows 37-46) to 工作表1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "公棕呺",
    "公棕呺 大神探路",
    "关注工众号",
    "关注公棕呺",
    "叉叉695叉叉",
    "叉叉511叉",
    "叉511叉叉",
    "叉65叉45",
    "叉叉396叉",
    "叉65叉45叉"
)

$startRow = 37
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
    $ws.Cells.Item($row, 2).Value = 1
}

# Scroll the view down to roughly show the newly added rows, and
# leave the final selection on the last appended cell (B46).
$win = $excel.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
$ws.Range("B46").Select() | Out-Null
